# Qualifier 2 DC vs KKR.
# Enter this match's predictions/points into Sheet1 row 71 and the
# corresponding "Qualifier 2" coin figures in the settlement table
# (rows 82-90, column G), plus the Rank 1 / Rank 2 predictors on Sheet2.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# --- Sheet1, row 71: new match "DC vs KKR" ---
$ws1.Range("C71").Value = "DC vs KKR"

$ws1.Range("E71").Value = 40
$ws1.Range("H71").Value = 20
$ws1.Range("K71").Value = 80
$ws1.Range("N71").Value = 0
$ws1.Range("Q71").Value = 30
$ws1.Range("T71").Value = 100
$ws1.Range("W71").Value = 70
$ws1.Range("Z71").Value = 60
$ws1.Range("AC71").Value = 50

# --- Sheet1, settlement table rows 82-90: Qualifier 2 coins (column G) ---
$ws1.Range("G82").Value = 0
$ws1.Range("G83").Value = 5
$ws1.Range("G84").Value = 3
$ws1.Range("G85").Value = 0
$ws1.Range("G86").Value = 3
$ws1.Range("G87").Value = 0
$ws1.Range("G88").Value = 0
$ws1.Range("G89").Value = 0
$ws1.Range("G90").Value = 3

# --- Sheet1, row 91: total row now hard-pinned to 0 ---
$ws1.Range("U91").Formula = "=0"

# --- Sheet2, row 50/51: record Rank 1 / Rank 2 predictors for Qualifier 2 ---
$ws2.Range("J50").Value = "Balaji"
$ws2.Range("J51").Value = "Ram"

# --- Sheet2 view: active cell (leave Sheet2's scroll position untouched) ---
$ws2.Range("F56:G56").Select()

# --- Sheet1 becomes the active sheet / tab; AD75 is the selected cell ---
$ws1.Activate()
$ws1.Range("AD75").Select()

# --- Workbook window position ---
$excel.ActiveWindow.Left = 3220
$excel.ActiveWindow.Top = 1360
